{"js": "// Word Office.js (Word JavaScript API) edit script.\n// Applies the textual corrections described in the commit\n// \"styling changes/ no more registration\".\n//\n// Strategy: use Body.search() to locate each exact (whole) phrase and\n// Range.insertText(..., Word.InsertLocation.replace) to swap in the new\n// text. Using full-phrase search/replace is robust against the document's\n// existing run-splitting (e.g. runs split around <w:proofErr/> tags) and\n// lets the host reflow/merge runs naturally, same as typing a correction\n// in the Word UI would.\n\nasync function replaceOnce(context, searchText, newText, options) {\n  const body = context.document.body;\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    return 0;\n  }\n  // Replace only the first match unless caller explicitly wants \"all\".\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n  return 1;\n}\n\nasync function replaceAll(context, searchText, newText, options) {\n  const body = context.document.body;\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return results.items.length;\n}\n\n// 1. \"Inloggevens voor de dashboard:\" -> \"Inloggegevens voor het dashboard:\"\nawait replaceOnce(\n  context,\n  \"Inloggevens voor de dashboard:\",\n  \"Inloggegevens voor het dashboard:\"\n);\n\n// 2. Capitalize \"ga naar de applicatie...\" -> \"Ga naar de applicatie...\"\nawait replaceOnce(\n  context,\n  \"ga naar de applicatie waar je kan aangeven hoe je dag was geweest.\",\n  \"Ga naar de applicatie waar je kan aangeven hoe je dag was geweest.\"\n);\n\n// 3. Capitalize the lower-case table header \"opmerkingen\" -> \"Opmerkingen\"\n//    (there are 4 lower-case occurrences; the 5th one is already capitalized\n//    and must stay untouched, hence matchCase: true).\nawait replaceAll(context, \"opmerkingen\", \"Opmerkingen\", { matchCase: true });\n\n// 4. Fix typo \"huistijl\" -> \"huisstijl\"\nawait replaceOnce(\n  context,\n  \"Is de huistijl van de app blauw, geel en zwart\",\n  \"Is de huisstijl van de app blauw, geel en zwart\"\n);\n\n// 5. Capitalize \"nederlands\" -> \"Nederlands\"\nawait replaceOnce(\n  context,\n  \"Is de app in het nederlands\",\n  \"Is de app in het Nederlands\"\n);\n\n// 6. Fix \"smily's\" -> \"smileys\"\nawait replaceOnce(\n  context,\n  \"Kan je kiezen tussen drie verschillende smily\\u2019s \",\n  \"Kan je kiezen tussen drie verschillende smileys \"\n);\n\n// 7. Reword the login paragraph (de -> het, split into two sentences).\nawait replaceOnce(\n  context,\n  \"Ga naar de login pagina van de dashboard gebruik de inlog gegevens van de begin van de documentatie. \",\n  \"Ga naar de login pagina van het dashboard. Gebruik de inloggegevens van het begin van de documentatie. \"\n);\n\n// 8. \"Als je bent ingelog in de dashboard:\" -> \"Als je ingelogd bent in het dashboard:\"\nawait replaceOnce(\n  context,\n  \"Als je bent ingelog in de dashboard:\",\n  \"Als je ingelogd bent in het dashboard:\"\n);\n\n// 9. Fix typo \"ban\" -> \"van\"\nawait replaceOnce(\n  context,\n  \"Is de opbouw ban het menu logisch\",\n  \"Is de opbouw van het menu logisch\"\n);\n\n// 11. \"knop naast de mail address\" -> \"knop naast de mailadressen\"\nawait replaceOnce(\n  context,\n  \"Is er een knop naast de mail address\",\n  \"Is er een knop naast de mailadressen\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the textual corrections described in the commit\n# \"styling changes/ no more registration\".\n#\n# Strategy: use Range.Find/Execute with a Replacement.Text to locate each\n# exact (whole) phrase and swap in the corrected text - the Word UI\n# equivalent of Ctrl+H \"Replace All\" for a specific phrase. This is robust\n# against the document's existing run-splitting (e.g. runs split around\n# proofing-error marks) since Word reflows/merges the runs for us.\n#\n# NOTE: this interpreter only binds function parameters positionally, so\n# the helper below is called as Replace-Phrase <find> <replace> <matchCase>\n# (no -Find/-Replace named arguments).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Phrase($Find, $Replace, $MatchCase) {\n    $find = $d.Content.Find\n    $find.Text = $Find\n    $find.Replacement.Text = $Replace\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $MatchCase\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Replace: 2 = wdReplaceAll\n    $find.Execute($Find, $MatchCase, $false, $false, $false, $false, $true, 1, $false, $Replace, 2) | Out-Null\n}\n\n# 1. \"Inloggevens voor de dashboard:\" -> \"Inloggegevens voor het dashboard:\"\nReplace-Phrase \"Inloggevens voor de dashboard:\" \"Inloggegevens voor het dashboard:\" $true\n\n# 2. Capitalize \"ga naar de applicatie...\" -> \"Ga naar de applicatie...\"\nReplace-Phrase \"ga naar de applicatie waar je kan aangeven hoe je dag was geweest.\" \"Ga naar de applicatie waar je kan aangeven hoe je dag was geweest.\" $true\n\n# 3. Capitalize the lower-case table header \"opmerkingen\" -> \"Opmerkingen\"\n#    (there are 4 lower-case occurrences; the 5th one is already capitalized\n#    and must stay untouched, hence MatchCase = $true and wdReplaceAll).\nReplace-Phrase \"opmerkingen\" \"Opmerkingen\" $true\n\n# 4. Fix typo \"huistijl\" -> \"huisstijl\"\nReplace-Phrase \"Is de huistijl van de app blauw, geel en zwart\" \"Is de huisstijl van de app blauw, geel en zwart\" $true\n\n# 5. Capitalize \"nederlands\" -> \"Nederlands\"\nReplace-Phrase \"Is de app in het nederlands\" \"Is de app in het Nederlands\" $true\n\n# 6. Fix \"smily\u2019s\" -> \"smileys\"\nReplace-Phrase \"Kan je kiezen tussen drie verschillende smily\u2019s \" \"Kan je kiezen tussen drie verschillende smileys \" $true\n\n# 7. Reword the login paragraph (de -> het, split into two sentences).\nReplace-Phrase \"Ga naar de login pagina van de dashboard gebruik de inlog gegevens van de begin van de documentatie. \" \"Ga naar de login pagina van het dashboard. Gebruik de inloggegevens van het begin van de documentatie. \" $true\n\n# 8. \"Als je bent ingelog in de dashboard:\" -> \"Als je ingelogd bent in het dashboard:\"\nReplace-Phrase \"Als je bent ingelog in de dashboard:\" \"Als je ingelogd bent in het dashboard:\" $true\n\n# 9. Fix typo \"ban\" -> \"van\"\nReplace-Phrase \"Is de opbouw ban het menu logisch\" \"Is de opbouw van het menu logisch\" $true\n\n# 11. \"knop naast de mail address\" -> \"knop naast de mailadressen\"\nReplace-Phrase \"Is er een knop naast de mail address\" \"Is er een knop naast de mailadressen\" $true\n"}
